$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so purely-numeric-looking
# strings (e.g. "379.64") are not auto-converted to numbers by Excel,
# matching the original inlineStr cell type used in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '51.343.93'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '2.922.46'
$ws.Range("E3").Value = '  -2.46%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '379.64'
$ws.Range("E5").Value = '  +7.30%  '
$ws.Range("D6").Value = '102.35'
$ws.Range("E6").Value = '  -4.84%  '
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  -3.02%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.586'
$ws.Range("E9").Value = '  -4.15%  '
$ws.Range("D10").Value = '37.14'
$ws.Range("E10").Value = '  -2.86%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '0.0834'
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '3.394.53'
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '18.25'
$ws.Range("E14").Value = '  -4.46%  '
$ws.Range("D15").Value = '7.34'
$ws.Range("E15").Value = '  -3.47%  '
$ws.Range("D16").Value = '2.933.71'
$ws.Range("E16").Value = '  -2.01%  '
$ws.Range("D17").Value = '0.926'
$ws.Range("E17").Value = '  -9.18%  '
$ws.Range("D18").Value = '51.316.12'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").Value = '3.43'
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").Value = '7.33'
$ws.Range("E20").Value = '  -2.19%  '
$ws.Range("D21").Value = '12.91'
$ws.Range("E21").Value = '  -4.16%  '
$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("D23").Value = '68.32'
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("D24").Value = '261.30'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").Value = '0.171'
$ws.Range("E26").Value = '  -3.66%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '4.11'
$ws.Range("E27").Value = '  -3.99%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '25.62'
$ws.Range("E29").Value = '  -4.48%  '
$ws.Range("D30").Value = '7.13'
$ws.Range("E30").Value = '  -4.96%  '
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Value = '6.82'
$ws.Range("E31").Value = '  +4.51%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.101'
$ws.Range("E32").Value = '  -3.90%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '9.80'
$ws.Range("E33").Value = '  -4.48%  '
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").Value = '2.11'
$ws.Range("E34").Value = '  -3.25%  '
$ws.Range("D35").Value = '51.50'
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '34.14'
$ws.Range("E36").Value = '  -5.36%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0420'
$ws.Range("E38").Value = '  -4.03%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").Value = '  -8.92%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '17.01'
$ws.Range("E40").Value = '  -3.63%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '2.56'
$ws.Range("E41").Value = '  -7.32%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '1.82'
$ws.Range("E42").Value = '  -8.18%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.114'
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '122.12'
$ws.Range("E44").Value = '  -1.59%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '21.79'
$ws.Range("E45").Value = '  -5.20%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").Value = '0.269'
$ws.Range("E47").Value = '  +10.74%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.030.63'
$ws.Range("E48").Value = '  -4.25%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '2.31'
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '3.15'
$ws.Range("E50").Value = '  -5.57%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '3.224.97'
$ws.Range("E51").Value = '  -2.07%  '

# Restore default (style-less) formatting now that the text values are set,
# so the cells end up with no explicit style index, same as before the edit.
$ws.Range("D2:D51").Style = "Normal"
